$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 1.09
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 1.5
$ws.Range("I4").Value = 1.71
$ws.Range("J4").Value = 3.75
$ws.Range("K4").Value = 1000
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 1.58

$ws.Range("F5").Value = 2.7
$ws.Range("G5").Value = 4.3
$ws.Range("H5").Value = 2.32
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 2.56
$ws.Range("K5").Value = 5.1
$ws.Range("P5").Value = 1.36
$ws.Range("Q5").Value = 2.32

$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 2.56
$ws.Range("I6").Value = 2.96
$ws.Range("J6").Value = 2.8
$ws.Range("K6").Value = 3.25
$ws.Range("P6").Value = 1.48
$ws.Range("Q6").Value = 2.42

$ws.Range("F8").Value = 1.89
$ws.Range("G8").Value = 2.48
$ws.Range("H8").Value = 1.67
$ws.Range("I8").Value = 7
$ws.Range("J8").Value = 3.05
$ws.Range("P8").Value = 1.74
$ws.Range("Q8").Value = 1.8

$ws.Range("H9").Value = 2.46
$ws.Range("P9").Value = 1.47
$ws.Range("Q9").Value = 2.72

$ws.Range("F10").Value = 2.02
$ws.Range("G10").Value = 2.62
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 4.7
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = 3.6
$ws.Range("P10").Value = 1.25
$ws.Range("Q10").Value = 1.01

$ws.Range("F11").Value = 1.49
$ws.Range("G11").Value = 1.58
$ws.Range("H11").Value = 7
$ws.Range("J11").Value = 4
$ws.Range("K11").Value = 6
$ws.Range("P11").Value = 2.08
$ws.Range("Q11").Value = 1.63

$ws.Range("P12").Value = 1.83
$ws.Range("Q12").Value = 2

$ws.Range("F13").Value = 1.63
$ws.Range("G13").Value = 1.72
$ws.Range("H13").Value = 5.9
$ws.Range("I13").Value = 6.6
$ws.Range("K13").Value = 4.4
$ws.Range("Q13").Value = 1.87

$ws.Range("F14").Value = 1.71
$ws.Range("G14").Value = 1.85
$ws.Range("H14").Value = 4.5
$ws.Range("I14").Value = 5.2
$ws.Range("Q14").Value = 1.61

$ws.Range("F15").Value = 4.8
$ws.Range("G15").Value = 5.5
$ws.Range("I15").Value = 2.02
$ws.Range("J15").Value = 3.3
$ws.Range("P15").Value = 1.56
$ws.Range("Q15").Value = 2.46

$ws.Range("F17").Value = 2.66
$ws.Range("G17").Value = 3.75
$ws.Range("H17").Value = 2.46
$ws.Range("I17").Value = 3.35
$ws.Range("J17").Value = 2.78
$ws.Range("K17").Value = 3.3
$ws.Range("P17").Value = 1.24

$ws.Range("F18").Value = 3.05
$ws.Range("G18").Value = 4.5
$ws.Range("H18").Value = 2.26
$ws.Range("I18").Value = 2.62
$ws.Range("J18").Value = 2.88
$ws.Range("K18").Value = 3.4
$ws.Range("P18").Value = 1.51
$ws.Range("Q18").Value = 2.32

$ws.Range("F19").Value = 1.91
$ws.Range("J19").Value = 3.75
$ws.Range("K19").Value = 4.9
$ws.Range("P19").Value = 2.9
$ws.Range("Q19").Value = 1.44

$ws.Range("F21").Value = 3.75
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 2.4
$ws.Range("I21").Value = 2.46
$ws.Range("J21").Value = 3
$ws.Range("K21").Value = 3.2
$ws.Range("P21").Value = 1.54
$ws.Range("Q21").Value = 2.68

$ws.Range("F22").Value = 2.06
$ws.Range("I22").Value = 4.7
$ws.Range("J22").Value = 3.35
$ws.Range("K22").Value = 4.3
$ws.Range("P22").Value = 1.79
$ws.Range("Q22").Value = 1.98

$ws.Range("F23").Value = 3.95
$ws.Range("G23").Value = 5.3
$ws.Range("H23").Value = 1.94
$ws.Range("I23").Value = 2.16
$ws.Range("J23").Value = 3.3
$ws.Range("K23").Value = 3.9
$ws.Range("P23").Value = 1.76
$ws.Range("Q23").Value = 2.06

$ws.Range("F24").Value = 1.25
$ws.Range("G24").Value = 13
$ws.Range("I24").Value = 1.5
$ws.Range("J24").Value = 4.4
$ws.Range("P24").Value = 1.83
$ws.Range("Q24").Value = 1.98

$ws.Range("H25").Value = 2.68
$ws.Range("J25").Value = 3.45
$ws.Range("P25").Value = 1.61

$ws.Range("H26").Value = 19
$ws.Range("J26").Value = 6.8
$ws.Range("K26").Value = 7
$ws.Range("P26").Value = 2.02
$ws.Range("Q26").Value = 1.82

$ws.Range("G29").Value = 7.6
$ws.Range("H29").Value = 1.66
$ws.Range("I29").Value = 1.8
$ws.Range("K29").Value = 4.2
$ws.Range("P29").Value = 1.75
$ws.Range("Q29").Value = 2.06

$ws.Range("F30").Value = 1.71
$ws.Range("G30").Value = 1.93
$ws.Range("H30").Value = 3.65
$ws.Range("I30").Value = 7
$ws.Range("K30").Value = 4.7
$ws.Range("Q30").Value = 1.84

$ws.Range("F31").Value = 2.3
$ws.Range("G31").Value = 2.64
$ws.Range("H31").Value = 2.84
$ws.Range("I31").Value = 4
$ws.Range("J31").Value = 3.2
$ws.Range("K31").Value = 3.7
$ws.Range("P31").Value = 1.66
$ws.Range("Q31").Value = 2.02

$ws.Range("F32").Value = 2.3
$ws.Range("G32").Value = 2.76
$ws.Range("I32").Value = 3.7
$ws.Range("J32").Value = 3.4
$ws.Range("P32").Value = 1.93
$ws.Range("Q32").Value = 1.74

$ws.Range("G33").Value = 2.26
$ws.Range("H33").Value = 3.25
$ws.Range("J33").Value = 3.25
$ws.Range("P33").Value = 1.85
$ws.Range("Q33").Value = 1.92

$ws.Range("F34").Value = 1.41
$ws.Range("H34").Value = 7.8
$ws.Range("I34").Value = 9.6
$ws.Range("P34").Value = 2.4
$ws.Range("Q34").Value = 1.56

$ws.Range("F35").Value = 1.66
$ws.Range("G35").Value = 1.93
$ws.Range("H35").Value = 3.7
$ws.Range("I35").Value = 6.6
$ws.Range("P35").Value = 2.1
$ws.Range("Q35").Value = 1.63

$ws.Range("A36").Value = "Uruguayan Primera Division"
$ws.Range("D36").Value = "Cerro"
$ws.Range("E36").Value = "Defensor Sporting"
$ws.Range("F36").Value = 3.5
$ws.Range("G36").Value = 5.1
$ws.Range("H36").Value = 2.04
$ws.Range("I36").Value = 2.62
$ws.Range("J36").Value = 2.78
$ws.Range("K36").Value = 3.35
$ws.Range("P36").Value = 1.42
$ws.Range("Q36").Value = 2.42

# Remove the last two rows (matches no longer tracked)
$ws.Range("A37:A38").EntireRow.Delete()

Write-Output "done"